$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.049.74'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.901.86'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9994'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.8376'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.20'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.65%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3293'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +3.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.66'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07062'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08090'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7653'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.919.68'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.271'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.49'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.032.86'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.15'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.871'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.46'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007778'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9995'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.152.22'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9996'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.003'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1744'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +23.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.291'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.58'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.40%  '
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.094'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.364'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.515'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05929'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +7.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.295'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.075'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.270'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7328'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.716'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01922'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.779'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4452'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.95'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.948'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8644'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.910'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9990'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.99'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.565'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.83%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.007.04'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.75%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.803'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.059.66'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.521'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.49%  '
